$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Roadmap")

$ws.Range("A8").Value = "Sean Duan"
$ws.Range("B8").Value = "Renderer"
$ws.Range("C8").Value = "Forward Pass实现基础的blin phong光照"

$ws.Range("E8").Select()
